$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Update existing D-column values (farms_total_count / farms_to_examine_count
# counts revised downward, plus the last existing week's remaining three
# variables and the trailing positive/negative counts) ---
$ws.Cells.Item(2, 4).Value = 11723
$ws.Cells.Item(3, 4).Value = 11468
$ws.Cells.Item(7, 4).Value = 11827
$ws.Cells.Item(8, 4).Value = 11114
$ws.Cells.Item(12, 4).Value = 11935
$ws.Cells.Item(13, 4).Value = 10612
$ws.Cells.Item(17, 4).Value = 11979
$ws.Cells.Item(18, 4).Value = 10061
$ws.Cells.Item(22, 4).Value = 12037
$ws.Cells.Item(23, 4).Value = 9633
$ws.Cells.Item(27, 4).Value = 12078
$ws.Cells.Item(28, 4).Value = 9173
$ws.Cells.Item(32, 4).Value = 12118
$ws.Cells.Item(33, 4).Value = 8765
$ws.Cells.Item(37, 4).Value = 12156
$ws.Cells.Item(38, 4).Value = 8241
$ws.Cells.Item(42, 4).Value = 12191
$ws.Cells.Item(43, 4).Value = 7674
$ws.Cells.Item(47, 4).Value = 12216
$ws.Cells.Item(48, 4).Value = 7126
$ws.Cells.Item(52, 4).Value = 12244
$ws.Cells.Item(53, 4).Value = 6484
$ws.Cells.Item(57, 4).Value = 12267
$ws.Cells.Item(58, 4).Value = 5873
$ws.Cells.Item(62, 4).Value = 12278
$ws.Cells.Item(63, 4).Value = 5707
$ws.Cells.Item(67, 4).Value = 12300
$ws.Cells.Item(68, 4).Value = 5523
$ws.Cells.Item(72, 4).Value = 12316
$ws.Cells.Item(73, 4).Value = 5087
$ws.Cells.Item(77, 4).Value = 12338
$ws.Cells.Item(78, 4).Value = 4609
$ws.Cells.Item(82, 4).Value = 12356
$ws.Cells.Item(83, 4).Value = 4090
$ws.Cells.Item(87, 4).Value = 12383
$ws.Cells.Item(88, 4).Value = 3684
$ws.Cells.Item(92, 4).Value = 12398
$ws.Cells.Item(93, 4).Value = 3345
$ws.Cells.Item(97, 4).Value = 12419
$ws.Cells.Item(98, 4).Value = 3034
$ws.Cells.Item(102, 4).Value = 12434
$ws.Cells.Item(103, 4).Value = 2758
$ws.Cells.Item(107, 4).Value = 12444
$ws.Cells.Item(108, 4).Value = 2471
$ws.Cells.Item(112, 4).Value = 12465
$ws.Cells.Item(113, 4).Value = 2234
$ws.Cells.Item(117, 4).Value = 12486
$ws.Cells.Item(118, 4).Value = 1989
$ws.Cells.Item(122, 4).Value = 12503
$ws.Cells.Item(123, 4).Value = 1735
$ws.Cells.Item(127, 4).Value = 12537
$ws.Cells.Item(128, 4).Value = 1476
$ws.Cells.Item(129, 4).Value = 11061
$ws.Cells.Item(130, 4).Value = 1472

# --- Append new week (202514 / LastDayOfWeek 2025-04-06) rows 132-136 ---
$newRows = @(
    @(202514, 45753, "farms_total_count", 12559),
    @(202514, 45753, "farms_to_examine_count", 1289),
    @(202514, 45753, "farms_examined_count", 11270),
    @(202514, 45753, "farms_examined_positive_count", 1452),
    @(202514, 45753, "farms_examined_negative_count", 9818)
)

$startRow = 132
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    # Reuse the date-formatted style already used by column B (copy format
    # from the last existing data row instead of assigning a fresh
    # NumberFormat, which would create a brand-new style entry).
    $ws.Cells.Item(131, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
$excel.CutCopyMode = $false

# --- Keep the view on the new tail of the data, matching the source edit ---
$ws.Application.ActiveWindow.ScrollRow = 112
$ws.Range("F132").Select()
